$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cells for team record: Wins / Losses / Ties (columns AD, AE, AF)
# Copy formatting from the existing header cell (AC1) so the new headers match
# the bold / bordered / centered header style used by the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record values (Wins=69, Losses=93, Ties=0) for every data row (2-50)
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = 69
    $ws.Cells.Item($r, 31).Value = 93
    $ws.Cells.Item($r, 32).Value = 0
}
